$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update input values (formulas in H, I, K, M and the summary K9 will recalc automatically)
$ws.Range("E3").Value = 77210.37
$ws.Range("F3").Value = 660.788
$ws.Range("E4").Value = 175782.39999999999
$ws.Range("F4").Value = 1517.5229999999999
$ws.Range("E5").Value = 302766.5
$ws.Range("F5").Value = 2689.365

# Update the selected cell/range on the sheet (was L13, now I11)
$ws.Activate()
$ws.Range("I11").Select()
